$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "сандальки"
$ws.Range("E2").Value = 300
$ws.Range("D6").Select()
